$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J
# Copy H1's formatting (bold, centered, bordered) onto I1 and J1,
# then overwrite the values with the new header text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows (2-18) - new values for columns I and J
$data = @{
    2  = @(8, 9)
    3  = @(8, 8)
    4  = @(5, 6)
    5  = @(5, 6)
    6  = @(1, 3)
    7  = @(6, 7)
    8  = @(5, 8)
    9  = @(1, 7)
    10 = @(1, 5)
    11 = @(1, 4)
    12 = @(1, 4)
    13 = @(1, 6)
    14 = @(1, 4)
    15 = @(1, 5)
    16 = @(1, 4)
    17 = @(1, 3)
    18 = @(4, 5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
